$d = $word.ActiveDocument

# The timestamp run "01:51 - 01:50***" must become "01:51 - 01:59***",
# with the final digit change ("0" -> "9") landing in its own run so the
# paragraph ends up as three runs:
#   "01:51 - 01:5"  +  "9"  +  "***"
# (all three runs keep identical Courier New rPr, matching the diff).

$target = "01:51 " + [char]0x2013 + " 01:50***"

$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute($target, $true, $false, $false, $false, $false,
                                  $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate target timestamp text"
}

# $findRange now spans the whole matched text "01:51 - 01:50***".
# The character to change is the "0" right before the trailing "***",
# i.e. the 4th-from-last character of the match.
$matchStart = $findRange.Start
$matchEnd = $findRange.End

$digitStart = $matchEnd - 4
$digitEnd = $digitStart + 1

$digitRange = $d.Range($digitStart, $digitEnd)

# Sanity check we grabbed the right character before mutating anything.
if ($digitRange.Text -ne "0") {
    throw "Unexpected character at split point: [$($digitRange.Text)]"
}

# Use a transient bookmark to force a hard run boundary around the single
# character we are about to replace, so it survives as its own <w:r> once
# the bookmark is removed (adjacent identical-formatting runs would
# otherwise be re-merged).
$bookmarkName = "TmpSplitMark"
$d.Bookmarks.Add($bookmarkName, $digitRange) | Out-Null

$digitRange.Text = "9"

$d.Bookmarks($bookmarkName).Delete()
